{"js": "// Widen the \"Name\" (surname) column by 2mm and narrow the \"Vorname\"\n// (first name) column by 2mm in the attendance-list table.\n// 1418 twips -> 1548 twips (Name)   i.e. 70.9pt -> 77.4pt\n// 1418 twips -> 1288 twips (Vorname) i.e. 70.9pt -> 64.4pt\n// (20 twips == 1 point; 1 mm == 56.6929... twips)\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// `TableCell.columnWidth` resizes the whole column (every row's cell in\n// that column plus the table's grid definition) in this host, matching\n// Word's own behaviour when you drag a column border \u2014 so touching row 0\n// is enough to update the entire column uniformly.\nconst nameCell = table.getCell(0, 2); // 3rd column: \"Name\" (surname)\nnameCell.columnWidth = 1548 / 20; // twips -> points\n\nconst firstNameCell = table.getCell(0, 3); // 4th column: \"Vorname\"\nfirstNameCell.columnWidth = 1288 / 20; // twips -> points\n\nawait context.sync();\n", "ps1": "# Widen the \"Name\" (surname) column by 2mm and narrow the \"Vorname\"\n# (first name) column by 2mm in the attendance-list table.\n# 1418 twips -> 1548 twips (Name)    i.e. 70.9pt -> 77.4pt\n# 1418 twips -> 1288 twips (Vorname) i.e. 70.9pt -> 64.4pt\n# (20 twips == 1 point)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Column.Width resizes the whole column (every row's cell in that column\n# plus the table's grid definition), matching Word's own behaviour when\n# dragging a column border.\n$nameColumn = $t.Columns.Item(3)       # 3rd column: \"Name\" (surname)\n$nameColumn.Width = 1548 / 20          # twips -> points\n\n$firstNameColumn = $t.Columns.Item(4)  # 4th column: \"Vorname\"\n$firstNameColumn.Width = 1288 / 20     # twips -> points\n"}
